# Morocco_FX.xlsx update:
#  - row 314 (2023-08-01): revise "high" (D) and "close" (F) values
#  - append three new monthly rows (315-317) for Sep/Oct/Nov 2023
#  - dimension grows from A1:G314 to A1:G317

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise existing row 314 ---
$ws.Range("D314").Value = 10.1674
$ws.Range("F314").Value = 10.1604

# --- Append row 315 (2023-09-01) ---
# Clone formatting from row 314 (date style, etc.) before writing values so the
# new row matches the existing table's look (s="2" date style on column A).
$ws.Range("A314:G314").Copy()
$ws.Range("A315:G315").PasteSpecial(-4122)
$ws.Range("A315").Value = 45170.33333333334
$ws.Range("B315").Value = "FX_IDC:USDMAD"
$ws.Range("C315").Value = 10.1598
$ws.Range("D315").Value = 10.3323
$ws.Range("E315").Value = 10.112
$ws.Range("F315").Value = 10.2833
$ws.Range("G315").Value = 0

# --- Append row 316 (2023-10-02) ---
$ws.Range("A314:G314").Copy()
$ws.Range("A316:G316").PasteSpecial(-4122)
$ws.Range("A316").Value = 45201.375
$ws.Range("B316").Value = "FX_IDC:USDMAD"
$ws.Range("C316").Value = 10.2833
$ws.Range("D316").Value = 10.3464
$ws.Range("E316").Value = 10.1978
$ws.Range("F316").Value = 10.2975
$ws.Range("G316").Value = 0

# --- Append row 317 (2023-11-01) ---
$ws.Range("A314:G314").Copy()
$ws.Range("A317:G317").PasteSpecial(-4122)
$ws.Range("A317").Value = 45231.375
$ws.Range("B317").Value = "FX_IDC:USDMAD"
$ws.Range("C317").Value = 10.2958
$ws.Range("D317").Value = 10.3303
$ws.Range("E317").Value = 10.1652
$ws.Range("F317").Value = 10.2043
$ws.Range("G317").Value = 0
